$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1, G1, H1 - set their text, then copy the formatting
# (bold, centered, bordered) from the existing header row (A1) so they match
# the look of the other header cells.
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Outlier flags (boolean) for rows 2-25 in columns F (KNN), G (SVM), H (RF).
# All FALSE except RF_Outliers_MAD for row 3 (Hb 3), which is TRUE.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    if ($r -eq 3) {
        $ws.Cells.Item($r, 8).Value = $true
    } else {
        $ws.Cells.Item($r, 8).Value = $false
    }
}
